$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'321.33"
$ws.Range("E2").Value = "'-3.49%"
$ws.Range("D3").Value = "'42.81"
$ws.Range("E3").Value = "'-6.03%"
$ws.Range("D4").Value = "'5.214"
$ws.Range("E4").Value = "'-4.90%"
$ws.Range("D5").Value = "'0.08234"
$ws.Range("E5").Value = "'-3.47%"
$ws.Range("D6").Value = "'4.328"
$ws.Range("E6").Value = "'-2.40%"
$ws.Range("D7").Value = "'1.777"
$ws.Range("E7").Value = "'-14.49%"
$ws.Range("D8").Value = "'0.9477"
$ws.Range("E8").Value = "'-4.26%"
$ws.Range("D9").Value = "'0.1122"
$ws.Range("E9").Value = "'-3.60%"
$ws.Range("D10").Value = "'0.1875"
$ws.Range("E10").Value = "'-2.21%"
$ws.Range("D11").Value = "'0.09410"
$ws.Range("E11").Value = "'-3.86%"
$ws.Range("E12").Value = "'-1.70%"
$ws.Range("D13").Value = "'7.455"
$ws.Range("E13").Value = "'-21.11%"
$ws.Range("E14").Value = "'-0.13%"
$ws.Range("D15").Value = "'0.001295"
$ws.Range("E15").Value = "'0.23%"
$ws.Range("D16").Value = "'0.005695"
$ws.Range("E16").Value = "'-3.73%"
$ws.Range("D17").Value = "'3.355"
$ws.Range("E17").Value = "'-1.07%"
$ws.Range("E18").Value = "'-0.47%"
$ws.Range("E19").Value = "'0.27%"
$ws.Range("D20").Value = "'0.1390"
$ws.Range("E20").Value = "'1.14%"
$ws.Range("D21").Value = "'0.2551"
$ws.Range("E21").Value = "'0.18%"
$ws.Range("D22").Value = "'0.04170"
$ws.Range("E22").Value = "'0.59%"
$ws.Range("E23").Value = "'-3.84%"
$ws.Range("D24").Value = "'0.004286"
$ws.Range("E24").Value = "'-3.98%"
$ws.Range("E25").Value = "'-6.09%"
$ws.Range("D26").Value = "'0.0002981"
$ws.Range("E26").Value = "'-0.21%"
$ws.Range("D38").Value = "'0.02661"
$ws.Range("E38").Value = "'-3.89%"
$ws.Range("D39").Value = "'0.05629"
$ws.Range("E39").Value = "'-1.97%"
$ws.Range("D40").Value = "'0.008147"
$ws.Range("E40").Value = "'3.96%"
$ws.Range("D41").Value = "'0.1401"
$ws.Range("E41").Value = "'-2.42%"
$ws.Range("D42").Value = "'0.006536"
$ws.Range("E42").Value = "'-9.85%"
$ws.Range("D43").Value = "'0.002119"
$ws.Range("E43").Value = "'0.55%"
$ws.Range("D44").Value = "'0.007694"
$ws.Range("E44").Value = "'-2.77%"
$ws.Range("D45").Value = "'0.3481"
$ws.Range("E45").Value = "'2.01%"
$ws.Range("D46").Value = "'0.00006763"
$ws.Range("E46").Value = "'-3.16%"
$ws.Range("E47").Value = "'0.06%"
$ws.Range("D48").Value = "'0.003075"
$ws.Range("E48").Value = "'-10.91%"
$ws.Range("D49").Value = "'0.004103"
$ws.Range("E49").Value = "'16.22%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'0.06%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.06%"
